$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), matching the style of H1 (bold/border header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I and J, rows 2-22
$data = @{
    2  = @(6, 7)
    3  = @(8, 8)
    4  = @(5, 6)
    5  = @(9, 9)
    6  = @(6, 7)
    7  = @(6, 8)
    8  = @(6, 7)
    9  = @(8, 8)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(8, 8)
    15 = @(8, 8)
    16 = @(4, 4)
    17 = @(7, 7)
    18 = @(9, 9)
    19 = @(9, 9)
    20 = @(7, 7)
    21 = @(7, 7)
    22 = @(3, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
